$wb = $excel.ActiveWorkbook

# --- summary_statistics ---
$ws = $wb.Worksheets.Item("summary_statistics")
$ws.Range("B2").Value = 323
$ws.Range("D2").Value = 18.93
$ws.Range("E2").Value = 1.35
$ws.Range("F2").Value = 0.73
$ws.Range("G2").Value = 2.27
$ws.Range("H2").Value = 2.34
$ws.Range("J2").Value = 2.340343958812575
$ws.Range("B3").Value = 323
$ws.Range("B4").Value = 323
$ws.Range("B5").Value = 323
$ws.Range("B6").Value = 323
$ws.Range("E6").Value = 0.13
$ws.Range("G6").Value = 0.34
$ws.Range("B7").Value = 323
$ws.Range("E7").Value = 0.39
$ws.Range("B8").Value = 323
$ws.Range("B9").Value = 323
$ws.Range("B10").Value = 323
$ws.Range("B11").Value = 323
$ws.Range("E11").Value = 0.05
$ws.Range("B12").Value = 323
$ws.Range("B13").Value = 323
$ws.Range("E13").Value = 818.99
$ws.Range("F13").Value = 189.05
$ws.Range("G13").Value = 2312.54
$ws.Range("H13").Value = 501.69
$ws.Range("I13").Value = 27.64283333333334
$ws.Range("J13").Value = 529.3333333333333
$ws.Range("B14").Value = 320
$ws.Range("E14").Value = 1.55
$ws.Range("G14").Value = 9.1
$ws.Range("B15").Value = 320
$ws.Range("E15").Value = 4.16
$ws.Range("F15").Value = 1.1
$ws.Range("G15").Value = 7.88
$ws.Range("B16").Value = 320
$ws.Range("D16").Value = 99.5
$ws.Range("E16").Value = 13.43
$ws.Range("G16").Value = 21.96
$ws.Range("H16").Value = 14.78
$ws.Range("J16").Value = 15.475
$ws.Range("B17").Value = 320
$ws.Range("E17").Value = 24.93
$ws.Range("F17").Value = 14.45
$ws.Range("G17").Value = 26.26
$ws.Range("H17").Value = 28.08
$ws.Range("J17").Value = 34.125
$ws.Range("B18").Value = 320
$ws.Range("E18").Value = 0.23
$ws.Range("G18").Value = 0.71
$ws.Range("B19").Value = 320
$ws.Range("B20").Value = 320
$ws.Range("E20").Value = 3.85
$ws.Range("F20").Value = 3.35
$ws.Range("H20").Value = 4.1
$ws.Range("I20").Value = 1.4
$ws.Range("B21").Value = 320
$ws.Range("E21").Value = 51.83
$ws.Range("F21").Value = 56.95
$ws.Range("G21").Value = 31.54
$ws.Range("H21").Value = 54.12
$ws.Range("I21").Value = 25.875
$ws.Range("J21").Value = 80
$ws.Range("B22").Value = 313
$ws.Range("E22").Value = 50.83
$ws.Range("F22").Value = 49.2
$ws.Range("G22").Value = 28.48
$ws.Range("H22").Value = 46.9
$ws.Range("I22").Value = 27.7
$ws.Range("J22").Value = 74.59999999999999
$ws.Range("B23").Value = 323
$ws.Range("E23").Value = 51.64
$ws.Range("G23").Value = 27.98
$ws.Range("H23").Value = 48.79
$ws.Range("J23").Value = 75.65000000000001

# --- dichotomous_stats ---
$ws = $wb.Worksheets.Item("dichotomous_stats")
$ws.Range("B2").Value = 285
$ws.Range("C2").Value = 38
$ws.Range("D2").Value = 0.607
$ws.Range("E2").Value = 1.421
$ws.Range("F2").Value = 0.8129999999999999
$ws.Range("G2").Value = 2.288
$ws.Range("H2").Value = 0.025
$ws.Range("I2").Value = 69.10299999999999
$ws.Range("J2").Value = 0.078
$ws.Range("K2").Value = 1.137
$ws.Range("B3").Value = 239
$ws.Range("C3").Value = 84
$ws.Range("D3").Value = 0.243
$ws.Range("E3").Value = 1.412
$ws.Range("F3").Value = 1.17
$ws.Range("G3").Value = 0.962
$ws.Range("H3").Value = 0.337
$ws.Range("I3").Value = 191.109
$ws.Range("K3").Value = 0.74
$ws.Range("B4").Value = 238
$ws.Range("C4").Value = 85
$ws.Range("D4").Value = 0.392
$ws.Range("E4").Value = 1.452
$ws.Range("F4").Value = 1.061
$ws.Range("G4").Value = 1.687
$ws.Range("H4").Value = 0.093
$ws.Range("I4").Value = 235.698
$ws.Range("J4").Value = -0.066
$ws.Range("K4").Value = 0.849
$ws.Range("B5").Value = 281
$ws.Range("C5").Value = 42
$ws.Range("D5").Value = 0.57
$ws.Range("E5").Value = 1.423
$ws.Range("F5").Value = 0.854
$ws.Range("G5").Value = 1.751
$ws.Range("H5").Value = 0.08500000000000001
$ws.Range("I5").Value = 60.423
$ws.Range("J5").Value = -0.081
$ws.Range("K5").Value = 1.22
$ws.Range("B6").Value = 196
$ws.Range("C6").Value = 127
$ws.Range("D6").Value = 0.264
$ws.Range("E6").Value = 1.453
$ws.Range("F6").Value = 1.189
$ws.Range("G6").Value = 1.104
$ws.Range("H6").Value = 0.27
$ws.Range("I6").Value = 319.582
$ws.Range("J6").Value = -0.207
$ws.Range("K6").Value = 0.735
$ws.Range("B7").Value = 230
$ws.Range("C7").Value = 93
$ws.Range("D7").Value = 0.317
$ws.Range("E7").Value = 1.44
$ws.Range("F7").Value = 1.124
$ws.Range("G7").Value = 1.325
$ws.Range("H7").Value = 0.187
$ws.Range("I7").Value = 243.955
$ws.Range("J7").Value = -0.154
$ws.Range("K7").Value = 0.787
$ws.Range("B8").Value = 220
$ws.Range("C8").Value = 103
$ws.Range("D8").Value = 0.301
$ws.Range("E8").Value = 1.445
$ws.Range("F8").Value = 1.144
$ws.Range("G8").Value = 1.263
$ws.Range("H8").Value = 0.207
$ws.Range("I8").Value = 275.203
$ws.Range("J8").Value = -0.168
$ws.Range("K8").Value = 0.769
$ws.Range("B9").Value = 244
$ws.Range("C9").Value = 79
$ws.Range("D9").Value = -0.049
$ws.Range("E9").Value = 1.337
$ws.Range("F9").Value = 1.386
$ws.Range("G9").Value = -0.206
$ws.Range("H9").Value = 0.837
$ws.Range("I9").Value = 207.32
$ws.Range("J9").Value = -0.515
$ws.Range("K9").Value = 0.418
$ws.Range("B10").Value = 308
$ws.Range("C10").Value = 15
$ws.Range("D10").Value = 0.581
$ws.Range("E10").Value = 1.376
$ws.Range("F10").Value = 0.796
$ws.Range("G10").Value = 1.308
$ws.Range("H10").Value = 0.208
$ws.Range("I10").Value = 16.767
$ws.Range("J10").Value = -0.357
$ws.Range("K10").Value = 1.518
$ws.Range("B11").Value = 223
$ws.Range("C11").Value = 100
$ws.Range("D11").Value = 0.529
$ws.Range("E11").Value = 1.513
$ws.Range("F11").Value = 0.984
$ws.Range("G11").Value = 2.344
$ws.Range("H11").Value = 0.02
$ws.Range("I11").Value = 296.367
$ws.Range("J11").Value = 0.08500000000000001
$ws.Range("K11").Value = 0.974

# --- anovas ---
$ws = $wb.Worksheets.Item("anovas")
$ws.Range("C2").Value = 154.4598779729278
$ws.Range("D2").Value = 51.48662599097594
$ws.Range("E2").Value = 10.9615182896583
$ws.Range("F2").Value = [double]"7.186355722274816e-07"
$ws.Range("C3").Value = 9.759004555692567
$ws.Range("D3").Value = 2.439751138923142
$ws.Range("E3").Value = 0.4721941557990477
$ws.Range("F3").Value = 0.7561471833188502
$ws.Range("C4").Value = 413.7763613511093
$ws.Range("D4").Value = 21.77770322900575
$ws.Range("E4").Value = 5.325621312230179
$ws.Range("F4").Value = [double]"3.772747684233216e-11"

# --- continuous_correlations ---
$ws = $wb.Worksheets.Item("continuous_correlations")
$ws.Range("B2").Value = 0.018
$ws.Range("C2").Value = 0.325
$ws.Range("D2").Value = 0.745
$ws.Range("E2").Value = 321
$ws.Range("F2").Value = -0.091
$ws.Range("G2").Value = 0.127
$ws.Range("C3").Value = -0.845
$ws.Range("D3").Value = 0.398
$ws.Range("E3").Value = 318
$ws.Range("F3").Value = -0.156
$ws.Range("G3").Value = 0.063
$ws.Range("B4").Value = 0.006
$ws.Range("C4").Value = 0.113
$ws.Range("D4").Value = 0.91
$ws.Range("E4").Value = 318
$ws.Range("F4").Value = -0.103
$ws.Range("G4").Value = 0.116
$ws.Range("B5").Value = -0.013
$ws.Range("C5").Value = -0.224
$ws.Range("D5").Value = 0.823
$ws.Range("E5").Value = 318
$ws.Range("F5").Value = -0.122
$ws.Range("G5").Value = 0.097
$ws.Range("B6").Value = 0.08799999999999999
$ws.Range("C6").Value = 1.569
$ws.Range("D6").Value = 0.118
$ws.Range("E6").Value = 318
$ws.Range("F6").Value = -0.022
$ws.Range("G6").Value = 0.195
$ws.Range("B7").Value = -0.093
$ws.Range("C7").Value = -1.668
$ws.Range("D7").Value = 0.096
$ws.Range("E7").Value = 318
$ws.Range("G7").Value = 0.017
$ws.Range("B8").Value = -0.078
$ws.Range("C8").Value = -1.387
$ws.Range("D8").Value = 0.167
$ws.Range("E8").Value = 318
$ws.Range("F8").Value = -0.186
$ws.Range("G8").Value = 0.032
$ws.Range("B9").Value = -0.039
$ws.Range("C9").Value = -0.7
$ws.Range("D9").Value = 0.484
$ws.Range("E9").Value = 318
$ws.Range("F9").Value = -0.148
$ws.Range("G9").Value = 0.07099999999999999
$ws.Range("B10").Value = -0.046
$ws.Range("C10").Value = -0.821
$ws.Range("D10").Value = 0.412
$ws.Range("E10").Value = 318
$ws.Range("F10").Value = -0.155
$ws.Range("G10").Value = 0.064
$ws.Range("B11").Value = 0.016
$ws.Range("C11").Value = 0.276
$ws.Range("D11").Value = 0.782
$ws.Range("E11").Value = 311
$ws.Range("F11").Value = -0.095
$ws.Range("G11").Value = 0.126
$ws.Range("B12").Value = 0.052
$ws.Range("C12").Value = 0.93
$ws.Range("D12").Value = 0.353
$ws.Range("E12").Value = 321
$ws.Range("F12").Value = -0.058
$ws.Range("G12").Value = 0.16
